$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.329.99'
$ws.Range("E2").Value = '  -4.85%  '
$ws.Range("D3").Value = '3.257.12'
$ws.Range("E3").Value = '  -7.79%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''590.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.17%  '
$ws.Range("D6").Value = '''152.58'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -12.40%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.248.72'
$ws.Range("E8").Value = '  -7.93%  '
$ws.Range("E9").Value = '  -10.91%  '
$ws.Range("D10").Value = '''0.172'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -13.36%  '
$ws.Range("D11").Value = '''6.78'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.05%  '
$ws.Range("D12").Value = '''0.510'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -12.64%  '
$ws.Range("D13").Value = '''38.58'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -17.09%  '
$ws.Range("D14").Value = '''0.0000244'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -12.02%  '
$ws.Range("D15").Value = '3.780.24'
$ws.Range("E15").Value = '  -7.72%  '
$ws.Range("D16").Value = '67.420.15'
$ws.Range("E16").Value = '  -4.83%  '
$ws.Range("D17").Value = '''548.71'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -10.03%  '
$ws.Range("D18").Value = '3.263.86'
$ws.Range("E18").Value = '  -7.63%  '
$ws.Range("D19").Value = '''7.30'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -13.17%  '
$ws.Range("D20").Value = '''0.115'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.90%  '
$ws.Range("D21").Value = '''15.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -14.00%  '
$ws.Range("D22").Value = '''0.770'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -12.97%  '
$ws.Range("D23").Value = '''7.88'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -13.32%  '
$ws.Range("D24").Value = '''85.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -12.79%  '
$ws.Range("D25").Value = '''13.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -12.65%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").Value = '''3.23'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -14.28%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '''8.12'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -10.66%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '''29.57'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -12.46%  '
$ws.Range("E30").Value = '  -16.61%  '
$ws.Range("E31").Value = '  -10.74%  '
$ws.Range("D32").Value = '''1.15'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -12.14%  '
$ws.Range("D33").Value = '''548.14'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -13.54%  '
$ws.Range("D34").Value = '''6.69'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -17.66%  '
$ws.Range("D35").Value = '''5.77'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -15.19%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").Value = '''0.0450'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.23%  '
$ws.Range("D38").Value = '''53.59'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.03%  '
$ws.Range("D39").Value = '''0.0863'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -13.45%  '
$ws.Range("D40").Value = '''9.26'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -14.31%  '
$ws.Range("E41").Value = '  -12.13%  '
$ws.Range("D42").Value = '2.950.95'
$ws.Range("E42").Value = '  -12.05%  '
$ws.Range("E43").Value = '  -22.80%  '
$ws.Range("E44").Value = '  -15.12%  '
$ws.Range("D45").Value = '0.0₃0585'
$ws.Range("E45").Value = '  -20.04%  '
$ws.Range("D46").Value = '''26.71'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -16.51%  '
$ws.Range("E47").Value = '  -15.13%  '
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").Value = '''1.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.09%  '
$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").Value = '''2.36'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -20.60%  '
$ws.Range("D50").Value = '''127.57'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.37%  '
$ws.Range("E51").Value = '  -11.95%  '
